# Add support for delimited file export
# ------------------------------------------------------------
# The "Format" sheet gains a new header row for picking an export
# type (Fixed Length / Delimited) and, when delimited, a delimiter
# character chosen from a small lookup list placed off to the side
# in column R. The existing "Field Name / Field Length" table is
# pushed down by three rows to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Format")
$ws.Activate() | Out-Null

# Push the existing content (Field Name/Field Length header + the
# three data rows) down by three rows.
$ws.Rows("1:3").Insert()

# --- New row 1: "Export Type" / "Delimiter" headers (big bold font,
#     matching the style already used for the "Field Name" header).
$ws.Range("A1").Value = "Export Type"
$ws.Range("A1").Font.Size = 20
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Value = "Delimiter"
$ws.Range("B1").Font.Size = 20
$ws.Range("B1").Font.Bold = $true

# --- New row 2: the two dropdown cells.
$ws.Range("A2").Value = "Fixed Length"
$ws.Range("B2").Value = "'N/A"

# --- Lookup list of delimiter choices in column R, used as the
#     source range for the B2 dropdown.
$ws.Range("R1").Value = "N/A"
$ws.Range("R2").Value = ","
$ws.Range("R3").Value = ";"
$ws.Range("R4").Value = "!"
$ws.Range("R5").Value = "TAB"
$ws.Columns("Q").ColumnWidth = 17.666666666666668

# --- Data validation dropdowns.
$ws.Range("A2").Validation.Add(3, 1, 1, """Fixed Length,Delimited""")
$ws.Range("A2").Validation.IgnoreBlank = $false
$ws.Range("B2").Validation.Add(3, 1, 1, "=`$R`$1:`$R`$5")
$ws.Range("B2").Validation.IgnoreBlank = $false

# --- Keep the dimension/page setup metadata tidy, matching the
#     "Field Name" sheet's portrait page setup.
$ws.PageSetup.Orientation = 1

# --- Leave the cursor parked on the new dropdown cell, as captured.
$ws.Range("D2").Select() | Out-Null

Write-Host "Delimited export support added to Format sheet"
